$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: seed the shared-string table in the exact order the new phrases were
# originally authored, by writing them once to a scratch column (Z) far outside the
# used range. The workbook already references 1168 unique strings (indices 0-1167);
# writing here first registers indices 1168-1247 in that order.
$ws.Range("Z1").Value = 'Consegui resolver tudo antes do prazo e fiquei super aliviado.'
$ws.Range("Z2").Value = 'A torcida inteira vibrou comigo quando acertei o lance final.'
$ws.Range("Z3").Value = 'O clima lá em casa hoje está tão leve que dá vontade de ficar conversando por horas.'
$ws.Range("Z4").Value = 'Recebi uma mensagem que fez meu dia começar muito melhor. '
$ws.Range("Z5").Value = 'O professor elogiou meu esforço e isso me deixou radiante. '
$ws.Range("Z6").Value = 'Treinei bastante e finalmente meu corpo respondeu do jeito que eu queria. '
$ws.Range("Z7").Value = 'Minha equipe fechou o projeto com chave de ouro.'
$ws.Range("Z8").Value = 'Passei na consulta e o médico disse que está tudo ótimo comigo.'
$ws.Range("Z9").Value = 'A risada dos meus amigos me trouxe uma alegria enorme.'
$ws.Range("Z10").Value = 'Hoje a aula prática foi tão divertida que o tempo voou.'
$ws.Range("Z11").Value = 'Meu irmão me surpreendeu com algo que eu queria muito.'
$ws.Range("Z12").Value = 'As pequenas conquistas da semana estão me deixando orgulhoso.'
$ws.Range("Z13").Value = 'As crianças da minha família espalham alegria pela casa inteira.'
$ws.Range("Z14").Value = 'Hoje consegui dormir bem e acordei com a mente leve.'
$ws.Range("Z15").Value = 'Terminei a série e amei cada segundo da história.'
$ws.Range("Z16").Value = 'Hoje me senti realmente vivo, presente e feliz.'
$ws.Range("Z17").Value = 'A comida ficou tão boa que todo mundo pediu a receita.'
$ws.Range("Z18").Value = 'O abraço que recebi hoje valeu mais do que mil palavras.'
$ws.Range("Z19").Value = 'Finalmente organizei meu quarto e ficou do jeito que eu queria.'
$ws.Range("Z20").Value = 'Ver meu time marcar o gol no último minuto foi pura emoção boa.'
$ws.Range("Z21").Value = 'Ouvi um barulho forte no portão e meu corpo travou.'
$ws.Range("Z22").Value = 'Minha respiração ficou curta quando percebi que estava sozinho na rua.'
$ws.Range("Z23").Value = 'A sensação de estar sendo seguido me deixa em pânico.'
$ws.Range("Z24").Value = 'O corredor escuro da escola sempre me dá um arrepio estranho.'
$ws.Range("Z25").Value = 'Não quero abrir aquele resultado médico, estou nervoso demais.'
$ws.Range("Z26").Value = 'Minha mão está tremendo desde que ouvi aquela notícia.'
$ws.Range("Z27").Value = 'A torcida adversária começou a gritar e eu fiquei com receio de confusão.'
$ws.Range("Z28").Value = 'A casa ficou silenciosa demais de repente.'
$ws.Range("Z29").Value = 'Meu coração dispara toda vez que penso no que pode acontecer.'
$ws.Range("Z30").Value = 'O cachorro começou a rosnar para um ponto vazio e isso me assustou.'
$ws.Range("Z31").Value = 'Não gosto de dirigir à noite, parece tudo mais perigoso.'
$ws.Range("Z32").Value = 'Ao ver a discussão aumentar, fiquei com medo de alguém se machucar.'
$ws.Range("Z33").Value = 'O clima na empresa ficou tão tenso que estou com receio de demissões.'
$ws.Range("Z34").Value = 'Aquela mensagem inesperada me deixou inquieto.'
$ws.Range("Z35").Value = 'O barulho vindo do quintal me fez congelar por um segundo.'
$ws.Range("Z36").Value = 'Quando a luz apagou, meu corpo inteiro arrepiou.'
$ws.Range("Z37").Value = 'Algo ali parecia errado, como se fosse melhor voltar.'
$ws.Range("Z38").Value = 'Senti que minha voz não ia sair se eu tentasse pedir ajuda.'
$ws.Range("Z39").Value = 'O vento forte batendo na janela me deu um susto enorme.'
$ws.Range("Z40").Value = 'Fico com medo de decepcionar quem acredita em mim.'
$ws.Range("Z41").Value = 'Já estou cansado dessa falta de respeito diária.'
$ws.Range("Z42").Value = 'Falam demais e resolvem de menos, isso irrita qualquer um.'
$ws.Range("Z43").Value = 'Meu time perdeu por erro de arbitragem, é revoltante.'
$ws.Range("Z44").Value = 'Não suporto quando ignoram tudo o que eu digo.'
$ws.Range("Z45").Value = 'Essa situação absurda não deveria nem existir.'
$ws.Range("Z46").Value = 'A escola prometeu melhorias e não cumpriu nada.'
$ws.Range("Z47").Value = 'Estou com uma vontade enorme de mandar todo mundo parar com essa palhaçada.'
$ws.Range("Z48").Value = 'Se continuarem me empurrando responsabilidades, eu vou explodir.'
$ws.Range("Z49").Value = 'Não dá pra aceitar tanta injustiça de braços cruzados.'
$ws.Range("Z50").Value = 'Já cansei de carregar problemas que não são meus.'
$ws.Range("Z51").Value = 'É irritante ver gente fingindo que está tudo normal.'
$ws.Range("Z52").Value = 'A forma como me trataram hoje foi simplesmente inaceitável.'
$ws.Range("Z53").Value = 'Minha paciência acabou com esse pessoal irresponsável.'
$ws.Range("Z54").Value = 'A torcida rival provocando sem parar me deixou furioso.'
$ws.Range("Z55").Value = 'Ninguém assume os erros, jogam tudo pra cima de mim.'
$ws.Range("Z56").Value = 'Esse atraso ridículo acabou com toda minha organização.'
$ws.Range("Z57").Value = 'Se repetirem isso mais uma vez, eu não garanto nada.'
$ws.Range("Z58").Value = 'O absurdo que ouvi hoje ainda está entalado na minha garganta.'
$ws.Range("Z59").Value = 'O caos no trânsito está me tirando do sério.'
$ws.Range("Z60").Value = 'Detesto quando fazem pouco caso dos meus sentimentos.'
$ws.Range("Z61").Value = 'Eu tento seguir em frente, mas parece que algo sempre me puxa de volta.'
$ws.Range("Z62").Value = 'Os corredores da escola hoje pareciam mais vazios que o normal.'
$ws.Range("Z63").Value = 'Senti um aperto no peito ao lembrar de como as coisas eram antes.'
$ws.Range("Z64").Value = 'A conversa que tive com minha família ainda ecoa na minha mente.'
$ws.Range("Z65").Value = 'Não consegui comemorar a vitória do time, algo dentro de mim está pesado.'
$ws.Range("Z66").Value = 'Ver as pessoas ao meu redor rindo me fez me sentir deslocado.'
$ws.Range("Z67").Value = 'A consulta médica trouxe mais dúvidas do que respostas.'
$ws.Range("Z68").Value = 'O silêncio do meu quarto ficou mais alto hoje.'
$ws.Range("Z69").Value = 'Minha energia sumiu, como se eu estivesse carregando o mundo nas costas.'
$ws.Range("Z70").Value = 'As mensagens que não chegam doem mais do que as que chegam.'
$ws.Range("Z71").Value = 'Sinto saudade de tempos que nunca vão voltar.'
$ws.Range("Z72").Value = 'A reunião no trabalho me deixou com uma sensação amarga.'
$ws.Range("Z73").Value = 'Mesmo rodeado de gente, continuo me sentindo invisível.'
$ws.Range("Z74").Value = 'Hoje me olhei no espelho e não reconheci quem estava lá.'
$ws.Range("Z75").Value = 'Percebi que não estou conseguindo acompanhar o ritmo de ninguém.'
$ws.Range("Z76").Value = 'Algo dentro de mim quebrou e eu não sei como consertar.'
$ws.Range("Z77").Value = 'A lembrança daquele momento feliz voltou para me machucar.'
$ws.Range("Z78").Value = 'Sinto como se tudo estivesse desmoronando aos poucos.'
$ws.Range("Z79").Value = 'Olhei o campo vazio e lembrei de todas as derrotas.'
$ws.Range("Z80").Value = 'O dia terminou, mas a sensação ruim ficou comigo.'

# --- Step 2: remove the scratch helper column; the 80 strings remain in the shared-
# string table (still referenced below) but Z stops showing up as used cells.
$ws.Range("Z1:Z80").ClearContents()

# --- Step 3: append the 80 new "frase"/"sentimento" rows (1165-1244) using the
# shared strings registered above.
$ws.Range("A1165").Value = 'Consegui resolver tudo antes do prazo e fiquei super aliviado.'
$ws.Range("B1165").Value = 'alegria'
$ws.Range("A1166").Value = 'A torcida inteira vibrou comigo quando acertei o lance final.'
$ws.Range("B1166").Value = 'alegria'
$ws.Range("A1167").Value = 'O clima lá em casa hoje está tão leve que dá vontade de ficar conversando por horas.'
$ws.Range("B1167").Value = 'alegria'
$ws.Range("A1168").Value = 'Recebi uma mensagem que fez meu dia começar muito melhor. '
$ws.Range("B1168").Value = 'alegria'
$ws.Range("A1169").Value = 'O professor elogiou meu esforço e isso me deixou radiante. '
$ws.Range("B1169").Value = 'alegria'
$ws.Range("A1170").Value = 'Treinei bastante e finalmente meu corpo respondeu do jeito que eu queria. '
$ws.Range("B1170").Value = 'alegria'
$ws.Range("A1171").Value = 'Minha equipe fechou o projeto com chave de ouro.'
$ws.Range("B1171").Value = 'alegria'
$ws.Range("A1172").Value = 'Passei na consulta e o médico disse que está tudo ótimo comigo.'
$ws.Range("B1172").Value = 'alegria'
$ws.Range("A1173").Value = 'A risada dos meus amigos me trouxe uma alegria enorme.'
$ws.Range("B1173").Value = 'alegria'
$ws.Range("A1174").Value = 'Hoje a aula prática foi tão divertida que o tempo voou.'
$ws.Range("B1174").Value = 'alegria'
$ws.Range("A1175").Value = 'Ver meu time marcar o gol no último minuto foi pura emoção boa.'
$ws.Range("B1175").Value = 'alegria'
$ws.Range("A1176").Value = 'Finalmente organizei meu quarto e ficou do jeito que eu queria.'
$ws.Range("B1176").Value = 'alegria'
$ws.Range("A1177").Value = 'A comida ficou tão boa que todo mundo pediu a receita.'
$ws.Range("B1177").Value = 'alegria'
$ws.Range("A1178").Value = 'O abraço que recebi hoje valeu mais do que mil palavras.'
$ws.Range("B1178").Value = 'alegria'
$ws.Range("A1179").Value = 'Meu irmão me surpreendeu com algo que eu queria muito.'
$ws.Range("B1179").Value = 'alegria'
$ws.Range("A1180").Value = 'As pequenas conquistas da semana estão me deixando orgulhoso.'
$ws.Range("B1180").Value = 'alegria'
$ws.Range("A1181").Value = 'Hoje consegui dormir bem e acordei com a mente leve.'
$ws.Range("B1181").Value = 'alegria'
$ws.Range("A1182").Value = 'As crianças da minha família espalham alegria pela casa inteira.'
$ws.Range("B1182").Value = 'alegria'
$ws.Range("A1183").Value = 'Terminei a série e amei cada segundo da história.'
$ws.Range("B1183").Value = 'alegria'
$ws.Range("A1184").Value = 'Hoje me senti realmente vivo, presente e feliz.'
$ws.Range("B1184").Value = 'alegria'
$ws.Range("A1185").Value = 'Ouvi um barulho forte no portão e meu corpo travou.'
$ws.Range("B1185").Value = 'medo'
$ws.Range("A1186").Value = 'Minha respiração ficou curta quando percebi que estava sozinho na rua.'
$ws.Range("B1186").Value = 'medo'
$ws.Range("A1187").Value = 'A sensação de estar sendo seguido me deixa em pânico.'
$ws.Range("B1187").Value = 'medo'
$ws.Range("A1188").Value = 'O corredor escuro da escola sempre me dá um arrepio estranho.'
$ws.Range("B1188").Value = 'medo'
$ws.Range("A1189").Value = 'Não quero abrir aquele resultado médico, estou nervoso demais.'
$ws.Range("B1189").Value = 'medo'
$ws.Range("A1190").Value = 'Minha mão está tremendo desde que ouvi aquela notícia.'
$ws.Range("B1190").Value = 'medo'
$ws.Range("A1191").Value = 'A torcida adversária começou a gritar e eu fiquei com receio de confusão.'
$ws.Range("B1191").Value = 'medo'
$ws.Range("A1192").Value = 'A casa ficou silenciosa demais de repente.'
$ws.Range("B1192").Value = 'medo'
$ws.Range("A1193").Value = 'Meu coração dispara toda vez que penso no que pode acontecer.'
$ws.Range("B1193").Value = 'medo'
$ws.Range("A1194").Value = 'O cachorro começou a rosnar para um ponto vazio e isso me assustou.'
$ws.Range("B1194").Value = 'medo'
$ws.Range("A1195").Value = 'Não gosto de dirigir à noite, parece tudo mais perigoso.'
$ws.Range("B1195").Value = 'medo'
$ws.Range("A1196").Value = 'Ao ver a discussão aumentar, fiquei com medo de alguém se machucar.'
$ws.Range("B1196").Value = 'medo'
$ws.Range("A1197").Value = 'O clima na empresa ficou tão tenso que estou com receio de demissões.'
$ws.Range("B1197").Value = 'medo'
$ws.Range("A1198").Value = 'Aquela mensagem inesperada me deixou inquieto.'
$ws.Range("B1198").Value = 'medo'
$ws.Range("A1199").Value = 'O barulho vindo do quintal me fez congelar por um segundo.'
$ws.Range("B1199").Value = 'medo'
$ws.Range("A1200").Value = 'Quando a luz apagou, meu corpo inteiro arrepiou.'
$ws.Range("B1200").Value = 'medo'
$ws.Range("A1201").Value = 'Algo ali parecia errado, como se fosse melhor voltar.'
$ws.Range("B1201").Value = 'medo'
$ws.Range("A1202").Value = 'Senti que minha voz não ia sair se eu tentasse pedir ajuda.'
$ws.Range("B1202").Value = 'medo'
$ws.Range("A1203").Value = 'O vento forte batendo na janela me deu um susto enorme.'
$ws.Range("B1203").Value = 'medo'
$ws.Range("A1204").Value = 'Fico com medo de decepcionar quem acredita em mim.'
$ws.Range("B1204").Value = 'medo'
$ws.Range("A1205").Value = 'Já estou cansado dessa falta de respeito diária.'
$ws.Range("B1205").Value = 'raiva'
$ws.Range("A1206").Value = 'Falam demais e resolvem de menos, isso irrita qualquer um.'
$ws.Range("B1206").Value = 'raiva'
$ws.Range("A1207").Value = 'Meu time perdeu por erro de arbitragem, é revoltante.'
$ws.Range("B1207").Value = 'raiva'
$ws.Range("A1208").Value = 'Não suporto quando ignoram tudo o que eu digo.'
$ws.Range("B1208").Value = 'raiva'
$ws.Range("A1209").Value = 'Essa situação absurda não deveria nem existir.'
$ws.Range("B1209").Value = 'raiva'
$ws.Range("A1210").Value = 'A escola prometeu melhorias e não cumpriu nada.'
$ws.Range("B1210").Value = 'raiva'
$ws.Range("A1211").Value = 'Estou com uma vontade enorme de mandar todo mundo parar com essa palhaçada.'
$ws.Range("B1211").Value = 'raiva'
$ws.Range("A1212").Value = 'Se continuarem me empurrando responsabilidades, eu vou explodir.'
$ws.Range("B1212").Value = 'raiva'
$ws.Range("A1213").Value = 'Não dá pra aceitar tanta injustiça de braços cruzados.'
$ws.Range("B1213").Value = 'raiva'
$ws.Range("A1214").Value = 'Já cansei de carregar problemas que não são meus.'
$ws.Range("B1214").Value = 'raiva'
$ws.Range("A1215").Value = 'É irritante ver gente fingindo que está tudo normal.'
$ws.Range("B1215").Value = 'raiva'
$ws.Range("A1216").Value = 'A forma como me trataram hoje foi simplesmente inaceitável.'
$ws.Range("B1216").Value = 'raiva'
$ws.Range("A1217").Value = 'Minha paciência acabou com esse pessoal irresponsável.'
$ws.Range("B1217").Value = 'raiva'
$ws.Range("A1218").Value = 'A torcida rival provocando sem parar me deixou furioso.'
$ws.Range("B1218").Value = 'raiva'
$ws.Range("A1219").Value = 'Ninguém assume os erros, jogam tudo pra cima de mim.'
$ws.Range("B1219").Value = 'raiva'
$ws.Range("A1220").Value = 'Esse atraso ridículo acabou com toda minha organização.'
$ws.Range("B1220").Value = 'raiva'
$ws.Range("A1221").Value = 'Se repetirem isso mais uma vez, eu não garanto nada.'
$ws.Range("B1221").Value = 'raiva'
$ws.Range("A1222").Value = 'O absurdo que ouvi hoje ainda está entalado na minha garganta.'
$ws.Range("B1222").Value = 'raiva'
$ws.Range("A1223").Value = 'O caos no trânsito está me tirando do sério.'
$ws.Range("B1223").Value = 'raiva'
$ws.Range("A1224").Value = 'Detesto quando fazem pouco caso dos meus sentimentos.'
$ws.Range("B1224").Value = 'raiva'
$ws.Range("A1225").Value = 'Eu tento seguir em frente, mas parece que algo sempre me puxa de volta.'
$ws.Range("B1225").Value = 'tristeza'
$ws.Range("A1226").Value = 'Os corredores da escola hoje pareciam mais vazios que o normal.'
$ws.Range("B1226").Value = 'tristeza'
$ws.Range("A1227").Value = 'Senti um aperto no peito ao lembrar de como as coisas eram antes.'
$ws.Range("B1227").Value = 'tristeza'
$ws.Range("A1228").Value = 'A conversa que tive com minha família ainda ecoa na minha mente.'
$ws.Range("B1228").Value = 'tristeza'
$ws.Range("A1229").Value = 'Não consegui comemorar a vitória do time, algo dentro de mim está pesado.'
$ws.Range("B1229").Value = 'tristeza'
$ws.Range("A1230").Value = 'Ver as pessoas ao meu redor rindo me fez me sentir deslocado.'
$ws.Range("B1230").Value = 'tristeza'
$ws.Range("A1231").Value = 'A consulta médica trouxe mais dúvidas do que respostas.'
$ws.Range("B1231").Value = 'tristeza'
$ws.Range("A1232").Value = 'O silêncio do meu quarto ficou mais alto hoje.'
$ws.Range("B1232").Value = 'tristeza'
$ws.Range("A1233").Value = 'Minha energia sumiu, como se eu estivesse carregando o mundo nas costas.'
$ws.Range("B1233").Value = 'tristeza'
$ws.Range("A1234").Value = 'As mensagens que não chegam doem mais do que as que chegam.'
$ws.Range("B1234").Value = 'tristeza'
$ws.Range("A1235").Value = 'Sinto saudade de tempos que nunca vão voltar.'
$ws.Range("B1235").Value = 'tristeza'
$ws.Range("A1236").Value = 'A reunião no trabalho me deixou com uma sensação amarga.'
$ws.Range("B1236").Value = 'tristeza'
$ws.Range("A1237").Value = 'Mesmo rodeado de gente, continuo me sentindo invisível.'
$ws.Range("B1237").Value = 'tristeza'
$ws.Range("A1238").Value = 'Hoje me olhei no espelho e não reconheci quem estava lá.'
$ws.Range("B1238").Value = 'tristeza'
$ws.Range("A1239").Value = 'Percebi que não estou conseguindo acompanhar o ritmo de ninguém.'
$ws.Range("B1239").Value = 'tristeza'
$ws.Range("A1240").Value = 'Algo dentro de mim quebrou e eu não sei como consertar.'
$ws.Range("B1240").Value = 'tristeza'
$ws.Range("A1241").Value = 'A lembrança daquele momento feliz voltou para me machucar.'
$ws.Range("B1241").Value = 'tristeza'
$ws.Range("A1242").Value = 'Sinto como se tudo estivesse desmoronando aos poucos.'
$ws.Range("B1242").Value = 'tristeza'
$ws.Range("A1243").Value = 'Olhei o campo vazio e lembrei de todas as derrotas.'
$ws.Range("B1243").Value = 'tristeza'
$ws.Range("A1244").Value = 'O dia terminou, mas a sensação ruim ficou comigo.'
$ws.Range("B1244").Value = 'tristeza'

# --- Step 4: reproduce the stray underline formatting left on the (content-less)
# I1231 cell, picked up while the author was scrolling/clicking through the sheet.
$ws.Range("I1231").Font.Underline = $true

# --- Step 5: leave the selection on I1231, matching the final view state.
$ws.Range("I1231").Select() | Out-Null
